$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added at the top of the data block
# (row 223). Insert a blank row there, which shifts every existing
# record (rows 223-262) down by one, and populate the new row with the
# latest reading.
$ws.Rows("223:223").Insert()

$ws.Range("A223").Value = 8
$ws.Range("B223").Value = "Terminal La Palmera de La Serena"
$ws.Range("C223").Value = "Coquimbo"
$ws.Range("D223").Value = 44504
$ws.Range("E223").Value = 4
$ws.Range("F223").Value = 100114001
$ws.Range("G223").Value = "Papa"
$ws.Range("H223").Value = "Cardinal"
$ws.Range("I223").Value = "1a (cosecha)"
$ws.Range("J223").Value = 2400
$ws.Range("K223").Value = 12000
$ws.Range("L223").Value = 12500
$ws.Range("M223").Value = 12250
$ws.Range("N223").Value = "$/saco 25 kilos"
$ws.Range("O223").Value = "Provincia del Elquí"
$ws.Range("P223").Value = 490
$ws.Range("Q223").Value = 25
$ws.Range("R223").Value = "Hortaliza"
